$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("A5").Value = 1
$ws.Range("A9").Value = 2
$ws.Range("A11").Value = 2
$ws.Range("A14").Value = 1
$ws.Range("A16").Value = 2
$ws.Range("A18").Value = 1
$ws.Range("A21").Value = 1
$ws.Range("A22").Value = 2
$ws.Range("A25").Value = 1
$ws.Range("A26").Value = 1
$ws.Range("A28").Value = 2
$ws.Range("A29").Value = 1
$ws.Range("A32").Value = 1
$ws.Range("A33").Value = 1
$ws.Range("A35").Value = 2
$ws.Range("A38").Value = 1
$ws.Range("A45").Value = 2
$ws.Range("A48").Value = 1
$ws.Range("A49").Value = 2
$ws.Range("A50").Value = 1
$ws.Range("A53").Value = 2
$ws.Range("A55").Value = 1
$ws.Range("A56").Value = 2
$ws.Range("A58").Value = 2
$ws.Range("A59").Value = 2
$ws.Range("A65").Value = 2
$ws.Range("A67").Value = 1
$ws.Range("A69").Value = 1
$ws.Range("A70").Value = 2
$ws.Range("A72").Value = 2
$ws.Range("A73").Value = 1
$ws.Range("A74").Value = 1
$ws.Range("A75").Value = 1
$ws.Range("A76").Value = 1
$ws.Range("A84").Value = 1
$ws.Range("A85").Value = 1
$ws.Range("A86").Value = 1
$ws.Range("A88").Value = 2
$ws.Range("A92").Value = 1
$ws.Range("A93").Value = 1
$ws.Range("A94").Value = 2
$ws.Range("A96").Value = 2
$ws.Range("A101").Value = 1
$ws.Range("A102").Value = 2
$ws.Range("A105").Value = 1
$ws.Range("A106").Value = 2
$ws.Range("A108").Value = 1
$ws.Range("A109").Value = 1
$ws.Range("A110").Value = 1
$ws.Range("A112").Value = 2
$ws.Range("A113").Value = 2
$ws.Range("A114").Value = 2
$ws.Range("A117").Value = 1
$ws.Range("A124").Value = 1
$ws.Range("A125").Value = 1
$ws.Range("A127").Value = 2
$ws.Range("A129").Value = 2
$ws.Range("A133").Value = 2
$ws.Range("A134").Value = 1
$ws.Range("A137").Value = 2
$ws.Range("A138").Value = 1
$ws.Range("A140").Value = 1
$ws.Range("A141").Value = 1
$ws.Range("A142").Value = 2
$ws.Range("A146").Value = 2
$ws.Range("A147").Value = 2
$ws.Range("A149").Value = 2
$ws.Range("A150").Value = 1
$ws.Range("A151").Value = 1
$ws.Range("A160").Value = 2
$ws.Range("A162").Value = 1
$ws.Range("A168").Value = 2
$ws.Range("A170").Value = 2
$ws.Range("A172").Value = 2
$ws.Range("A173").Value = 2
$ws.Range("A174").Value = 1
$ws.Range("A175").Value = 1
$ws.Range("A177").Value = 2
$ws.Range("A186").Value = 2
$ws.Range("A188").Value = 1
$ws.Range("A190").Value = 2
$ws.Range("A194").Value = 1
$ws.Range("A197").Value = 1
$ws.Range("A199").Value = 1
$ws.Range("A200").Value = 2
